$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "max" column (C) entirely; this shifts "prediction" (D->C)
# and "rejection-f" (E->D) left by one column and updates the used
# dimension from A1:E2 to A1:D2.
$ws.Range("C1").EntireColumn.Delete()

# Update the numeric value in B2 (was the "1" count for a multi-child
# prediction; now reflects the single remaining child's score).
$ws.Range("B2").Value = 1.920968662069615
